$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I0 and IF columns, matching the style of the existing
# header cell (H1) — bold, centered, bordered.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-16: I column then J column values
$data = @(
    @(1, 5),
    @(1, 4),
    @(1, 3),
    @(9, 9),
    @(1, 6),
    @(1, 7),
    @(1, 6),
    @(1, 6),
    @(1, 4),
    @(1, 7),
    @(1, 7),
    @(8, 8),
    @(5, 7),
    @(5, 7),
    @(3, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
